$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("08-09-2021", 3768, 8140),
    @("09-09-2021", 3723, 8913),
    @("10-09-2021", 3727, 7243),
    @("13-09-2021", 3404, 6365),
    @("14-09-2021", 3626, 7789),
    @("15-09-2021", 3616, 10608)
)

$startRow = 174
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    # Enter the date-like text as a formula returning a string, then convert
    # it to a plain value in place. This avoids Excel's automatic
    # text-to-date conversion that would otherwise turn ambiguous strings
    # like "08-09-2021" into date serial numbers.
    $cellA.Formula = "=""" + $newData[$i][0] + """"
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}
